# Rename the worksheet from "Property1" to "DataNode" -- unifying the
# DataNode / DataTable / Entity naming used across the config workbooks.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Name = "DataNode"

# Leave the sheet with the same cursor position the author's Excel
# session recorded on save (bottom pane of the frozen-row split).
$ws.Range("L39").Select() | Out-Null
